# Property report workbook - issue #5: property boat&car done
# Adds the standard metadata columns (property_category, category, date,
# legislator_name, legislator_id, source_file, index) to the "汽車" (car)
# sheet, matching the schema already used on the 土地/建物 sheets, and
# relabels the previous "data-row-as-header" B1/C1/D1/E1 cells with the
# real column headers (name/capacity/owner/register_date/...).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("汽車")   # 汽車 (car) sheet

# ---- Header row (row 1) ----
$ws.Range("B1").Value = "name"
$ws.Range("C1").Value = "capacity"
$ws.Range("D1").Value = "owner"
$ws.Range("E1").Value = "register_date"
$ws.Range("F1").Value = "register_reason"
$ws.Range("G1").Value = "acquire_value"
$ws.Range("H1").Value = "property_category"
$ws.Range("I1").Value = "category"
$ws.Range("J1").Value = "date"
$ws.Range("K1").Value = "legislator_name"
$ws.Range("L1").Value = "legislator_id"
$ws.Range("M1").Value = "source_file"
$ws.Range("N1").Value = "index"

# Match the bold / centered / bordered header style used by columns B:G
$headerRange = $ws.Range("H1:N1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

# ---- Data row (row 2) ----
$ws.Range("B2").Value = "國瑞國產客貨車"
$ws.Range("E2").Value = "93年09月01日"
$ws.Range("H2").Value = "land"
$ws.Range("I2").Value = "normal"

# "date" must stay a plain text value ("2012-04-26"), not get auto-converted
# into a date serial number - force text format for the assignment, then
# drop the formatting override again so the cell matches its plain siblings.
$ws.Range("J2").NumberFormat = "@"
$ws.Range("J2").Value = "2012-04-26"
$ws.Range("J2").ClearFormats()

$ws.Range("K2").Value = "林鴻池"
$ws.Range("L2").Value = 1340
$ws.Range("M2").Value = "tmpdb4b1"
$ws.Range("N2").Value = 48
